$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 75.5
$ws.Range("I55").Value = 94
$ws.Range("J55").Value = 66.25
$ws.Range("K55").Value = 94
$ws.Range("L55").Value = 66.25
$ws.Range("M55").Value = 120
$ws.Range("N55").Value = -494.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 6640
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6640
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 19920
$ws.Range("N70").Value = -20460

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 6640
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6640
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 19920
$ws.Range("N73").Value = -21792

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -4594

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -3596

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 3639.2856
$ws.Range("I103").Value = 2331.2
$ws.Range("J103").Value = 4828.4546
$ws.Range("K103").Value = 6993.599999999999
$ws.Range("L103").Value = 14485.3638
$ws.Range("M103").Value = -6407.599999999999
$ws.Range("N103").Value = -15657.3638

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3999.6667
$ws.Range("I116").Value = 3999.5
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 3999.5
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -557.5
$ws.Range("N116").Value = -10884

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 13451
$ws.Range("I132").Value = 12576.333
$ws.Range("J132").Value = 16599.8
$ws.Range("K132").Value = 37728.999
$ws.Range("L132").Value = 49799.39999999999
$ws.Range("M132").Value = -35198.999
$ws.Range("N132").Value = -54859.39999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2744.375
$ws.Range("I138").Value = 2343.8
$ws.Range("J138").Value = 3412
$ws.Range("K138").Value = 7031.400000000001
$ws.Range("L138").Value = 10236
$ws.Range("M138").Value = -1891.400000000001
$ws.Range("N138").Value = -20516

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3241.7144
$ws.Range("I141").Value = 2023.75
$ws.Range("J141").Value = 4865.6665
$ws.Range("K141").Value = 6071.25
$ws.Range("L141").Value = 14596.9995
$ws.Range("M141").Value = -891.25
$ws.Range("N141").Value = -24956.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3175.762
$ws.Range("I45").Value = 2570.0588
$ws.Range("J45").Value = 5750
$ws.Range("K45").Value = 2570.0588
$ws.Range("L45").Value = 5750
$ws.Range("M45").Value = -2193.0588
$ws.Range("N45").Value = -6504

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3693.0667
$ws.Range("I61").Value = 3024.75
$ws.Range("J61").Value = 6366.3335
$ws.Range("K61").Value = 3024.75
$ws.Range("L61").Value = 6366.3335
$ws.Range("M61").Value = -2812.75
$ws.Range("N61").Value = -6790.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3377.3333
$ws.Range("I132").Value = 1679.4
$ws.Range("J132").Value = 5499.75
$ws.Range("K132").Value = 5038.200000000001
$ws.Range("L132").Value = 16499.25
$ws.Range("M132").Value = -2508.200000000001
$ws.Range("N132").Value = -21559.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3693.0667
$ws.Range("I136").Value = 3024.75
$ws.Range("J136").Value = 6366.3335
$ws.Range("K136").Value = 9074.25
$ws.Range("L136").Value = 19099.0005
$ws.Range("M136").Value = -6524.25
$ws.Range("N136").Value = -24199.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5218.778
$ws.Range("I86").Value = 4517.769
$ws.Range("J86").Value = 7041.4
$ws.Range("K86").Value = 4517.769
$ws.Range("L86").Value = 7041.4
$ws.Range("M86").Value = -3394.769
$ws.Range("N86").Value = -9287.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5218.778
$ws.Range("I89").Value = 4517.769
$ws.Range("J89").Value = 7041.4
$ws.Range("K89").Value = 22588.845
$ws.Range("L89").Value = 35207
$ws.Range("M89").Value = -16972.845
$ws.Range("N89").Value = -46439

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5515.25
$ws.Range("I134").Value = 4586
$ws.Range("J134").Value = 8303
$ws.Range("K134").Value = 13758
$ws.Range("L134").Value = 24909
$ws.Range("M134").Value = -11223
$ws.Range("N134").Value = -29979

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1392.5
$ws.Range("I17").Value = 471
$ws.Range("J17").Value = 6000
$ws.Range("K17").Value = 1413
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = -1244
$ws.Range("N17").Value = -18338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6255.6113
$ws.Range("I39").Value = 504.25
$ws.Range("J39").Value = 7898.857
$ws.Range("K39").Value = 1512.75
$ws.Range("L39").Value = 23696.571
$ws.Range("M39").Value = -1218.75
$ws.Range("N39").Value = -24284.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2665
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 2997.5
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 8992.5
$ws.Range("M63").Value = -5251
$ws.Range("N63").Value = -10490.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 2665
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 2997.5
$ws.Range("K66").Value = 18000
$ws.Range("L66").Value = 26977.5
$ws.Range("M66").Value = -14256
$ws.Range("N66").Value = -34465.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 1900000
$ws.Range("I128").Value = 1900000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 5700000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -5695020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2247.3
$ws.Range("I129").Value = 808
$ws.Range("J129").Value = 3686.6
$ws.Range("K129").Value = 2424
$ws.Range("L129").Value = 11059.8
$ws.Range("M129").Value = 2576
$ws.Range("N129").Value = -21059.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 9223
$ws.Range("I139").Value = 8953
$ws.Range("J139").Value = 10033
$ws.Range("K139").Value = 26859
$ws.Range("L139").Value = 30099
$ws.Range("M139").Value = -21719
$ws.Range("N139").Value = -40379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I80").Value = 2732.8333
$ws.Range("J80").Value = 2416.3333
$ws.Range("K80").Value = 2732.8333
$ws.Range("L80").Value = 2416.3333
$ws.Range("M80").Value = -1734.8333
$ws.Range("N80").Value = -4412.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I83").Value = 2732.8333
$ws.Range("J83").Value = 2416.3333
$ws.Range("K83").Value = 13664.1665
$ws.Range("L83").Value = 12081.6665
$ws.Range("M83").Value = -8672.166499999999
$ws.Range("N83").Value = -22065.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 733.6
$ws.Range("I97").Value = 812.2
$ws.Range("J97").Value = 419.2
$ws.Range("K97").Value = 812.2
$ws.Range("L97").Value = 419.2
$ws.Range("M97").Value = -316.2
$ws.Range("N97").Value = -1411.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1611.1471
$ws.Range("I102").Value = 1159.3667
$ws.Range("J102").Value = 4999.5
$ws.Range("K102").Value = 1159.3667
$ws.Range("L102").Value = 4999.5
$ws.Range("M102").Value = 462.6333
$ws.Range("N102").Value = -8243.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 100998
$ws.Range("I35").Value = 1247.5
$ws.Range("J35").Value = 500000
$ws.Range("K35").Value = 1247.5
$ws.Range("L35").Value = 500000
$ws.Range("M35").Value = -911.5
$ws.Range("N35").Value = -500672

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7431.125
$ws.Range("I46").Value = 9275
$ws.Range("J46").Value = 6816.5
$ws.Range("K46").Value = 9275
$ws.Range("L46").Value = 6816.5
$ws.Range("M46").Value = -9087
$ws.Range("N46").Value = -7192.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 657.913
$ws.Range("I55").Value = 1198.5
$ws.Range("J55").Value = 369.6
$ws.Range("K55").Value = 1198.5
$ws.Range("L55").Value = 369.6
$ws.Range("M55").Value = -1025.5
$ws.Range("N55").Value = -715.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2323.4666
$ws.Range("I93").Value = 2154.5386
$ws.Range("J93").Value = 3421.5
$ws.Range("K93").Value = 2154.5386
$ws.Range("L93").Value = 3421.5
$ws.Range("M93").Value = -906.5385999999999
$ws.Range("N93").Value = -5917.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4412.5
$ws.Range("I122").Value = 3897
$ws.Range("J122").Value = 5959
$ws.Range("K122").Value = 11691
$ws.Range("L122").Value = 17877
$ws.Range("M122").Value = -9241
$ws.Range("N122").Value = -22777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5332.5386
$ws.Range("I132").Value = 4943.1665
$ws.Range("J132").Value = 10005
$ws.Range("K132").Value = 14829.4995
$ws.Range("L132").Value = 30015
$ws.Range("M132").Value = -12299.4995
$ws.Range("N132").Value = -35075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 24384.666
$ws.Range("I32").Value = 7069.75
$ws.Range("J32").Value = 59014.5
$ws.Range("K32").Value = 7069.75
$ws.Range("L32").Value = 59014.5
$ws.Range("M32").Value = -6752.75
$ws.Range("N32").Value = -59648.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 25049.75
$ws.Range("I74").Value = 22099
$ws.Range("J74").Value = 26033.334
$ws.Range("K74").Value = 22099
$ws.Range("L74").Value = 26033.334
$ws.Range("M74").Value = -21163
$ws.Range("N74").Value = -27905.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 25049.75
$ws.Range("I77").Value = 22099
$ws.Range("J77").Value = 26033.334
$ws.Range("K77").Value = 66297
$ws.Range("L77").Value = 78100.00199999999
$ws.Range("M77").Value = -61617
$ws.Range("N77").Value = -87460.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2000
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -4696

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 560.8889
$ws.Range("I107").Value = 536.4666999999999
$ws.Range("J107").Value = 683
$ws.Range("K107").Value = 1609.4001
$ws.Range("L107").Value = 2049
$ws.Range("M107").Value = 310.5999000000002
$ws.Range("N107").Value = -5889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1349.92
$ws.Range("I122").Value = 1289
$ws.Range("J122").Value = 1669.75
$ws.Range("K122").Value = 3867
$ws.Range("L122").Value = 5009.25
$ws.Range("M122").Value = -1417
$ws.Range("N122").Value = -9909.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3609.5652
$ws.Range("I126").Value = 1740.7333
$ws.Range("J126").Value = 7113.625
$ws.Range("K126").Value = 5222.199900000001
$ws.Range("L126").Value = 21340.875
$ws.Range("M126").Value = -2752.199900000001
$ws.Range("N126").Value = -26280.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3350.6086
$ws.Range("I132").Value = 3139.318
$ws.Range("J132").Value = 7999
$ws.Range("K132").Value = 9417.954000000002
$ws.Range("L132").Value = 23997
$ws.Range("M132").Value = -6887.954000000002
$ws.Range("N132").Value = -29057
